$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "1" to "Mestia"
$ws.Name = "Mestia"

# The table used to have an empty row 8 between the data rows (5-7) and the
# footnote row (9). Remove that empty row so the footnote moves up to row 8.
$ws.Rows(8).Delete()

# Columns N (year 2022) on the "Urban" and "Rural" rows previously held
# explicit counts; they should show the confidential/unavailable marker
# like the rest of the row instead.
$ws.Range("N6").Value = "…"
$ws.Range("N7").Value = "…"
